$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B values (one per row, 1..50) taken from the source data.
$values = @(
    63,
    83,
    56,
    42,
    86,
    44,
    73,
    61.285714285714299,
    60.642857142857103,
    60,
    59.357142857142897,
    58.714285714285701,
    58.071428571428598,
    57.428571428571402,
    56.785714285714299,
    56.142857142857103,
    55.5,
    54.857142857142897,
    54.214285714285701,
    53.571428571428598,
    52.928571428571402,
    52.285714285714299,
    51.642857142857103,
    51,
    50.357142857142897,
    49.714285714285701,
    49.071428571428598,
    48.428571428571402,
    47.785714285714299,
    47.142857142857103,
    46.5,
    45.857142857142897,
    45.214285714285701,
    44.571428571428598,
    43.928571428571402,
    43.285714285714299,
    42.642857142857103,
    42,
    41.357142857142897,
    40.714285714285701,
    40.071428571428598,
    39.428571428571402,
    38.785714285714299,
    38.142857142857103,
    37.5,
    36.857142857142897,
    36.214285714285701,
    35.571428571428598,
    34.928571428571402,
    34.285714285714299
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the active selection to match the saved view state.
$ws.Range("C2").Select()
